{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (plus the blank separator paragraph right before it) that used to sit\n// right after the \"Requisitos\" entry (e.g. \"LOB1053: F\u00edsica III (Requisito\n// fraco)\"), while leaving the rest of the document (including the blank\n// paragraph + page-break paragraph that follow the footer) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the footer paragraphs by their exact text content.\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    jupiterIndex = i;\n  } else if (t.includes(\"Powered by Jekyll\")) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex !== -1 && copyrightIndex !== -1) {\n  // The blank separator paragraph immediately precedes the \"Ver no Jupiter\u2026\"\n  // paragraph; delete it too (but only if it really is empty, so we never\n  // eat an unrelated paragraph).\n  const candidateBlankIndex = jupiterIndex - 1;\n  const toDelete = [];\n  if (\n    candidateBlankIndex >= 0 &&\n    items[candidateBlankIndex].text.trim() === \"\"\n  ) {\n    toDelete.push(items[candidateBlankIndex]);\n  }\n  toDelete.push(items[jupiterIndex]);\n  toDelete.push(items[copyrightIndex]);\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" footer block\n# (plus the blank separator paragraph right before it) that used to sit\n# right after the \"Requisitos\" entry (e.g. \"LOB1053: Fisica III (Requisito\n# fraco)\"), while leaving the rest of the document (including the blank\n# paragraph + page-break paragraph that follow the footer) untouched.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$jupiterIdx = -1\n$copyrightIdx = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $jupiterIdx = $i\n    } elseif ($t.Contains(\"Powered by Jekyll\")) {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -ge 1 -and $copyrightIdx -ge 1) {\n    # Blank separator paragraph right before \"Ver no Jupiter...\", if present.\n    $blankIdx = $jupiterIdx - 1\n    $hasBlank = $false\n    if ($blankIdx -ge 1) {\n        $blankText = $d.Paragraphs.Item($blankIdx).Range.Text.Trim()\n        if ($blankText -eq \"\") {\n            $hasBlank = $true\n        }\n    }\n\n    # Delete from the bottom up so earlier indices stay valid.\n    $d.Paragraphs.Item($copyrightIdx).Range.Delete()\n    $d.Paragraphs.Item($jupiterIdx).Range.Delete()\n    if ($hasBlank) {\n        $d.Paragraphs.Item($blankIdx).Range.Delete()\n    }\n}\n"}
